$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column I: "Provide or Consume" header, with the EPG row's
# provide/consume value moved from the stray H2 cell into I2 as "consume".
$ws.Range("I1").Value = "Provide or Consume"
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = "consume"

# Matches the author's final selection in the saved workbook.
$ws.Range("G15").Select() | Out-Null
